$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (last-changed) date column C for rows 2-28: 45183 -> 45184
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Row 2 (record "A 30703-2023"): rewrite the document-link formulas in
# columns S, T, V, W, X, Y to add the display-text argument to HYPERLINK(...).
# Note: the S2 formula is reproduced exactly as authored, including its
# (malformed) unbalanced quoting.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/artfynd/A 30703-2023.xlsx, "A 30703-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/kartor/A 30703-2023.png", "A 30703-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/klagomål/A 30703-2023.docx", "A 30703-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/klagomålsmail/A 30703-2023.docx", "A 30703-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/tillsyn/A 30703-2023.docx", "A 30703-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/tillsynsmail/A 30703-2023.docx", "A 30703-2023")'
